$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the shared string used as the table name in C2: "ClienteTest" -> "CL_ClienteTest".
#    (All F-column formulas reference $C$2, so their cached values recalc automatically.)
$ws.Range("C2").Value = "CL_ClienteTest"

# 2) New header-ish cells on row 2: E2 = "User", F2 = "AUTO" (used as a literal by the new
#    G-column formulas below via the absolute reference $F$2).
$ws.Range("E2").Value = "User"
$ws.Range("F2").Value = "AUTO"

# 3) New column G: a PL/SQL-ish "INSERTA" call built from D (Nombre) and E (Direccion) plus
#    the literal in $F$2. Row 4 is the "anchor" formula (not shared); rows 5:36 are filled as
#    one range-assignment so the engine groups them into a single shared formula, matching how
#    Excel itself would store a fill-down.
$ws.Range("G4").Formula = '= "wcod_cliente_n := null; CL_PCLIENTETEST.INSERTA (wcod_cliente_n, ''" &D4 & "'',''" &E4 & "'', 2, NULL, ''" &  $F$2 &"'', NULL);"'
$ws.Range("G5:G36").Formula = '= "wcod_cliente_n := null; CL_PCLIENTETEST.INSERTA (wcod_cliente_n, ''" &D5 & "'',''" &E5 & "'', 2, NULL, ''" &  $F$2 &"'', NULL);"'

# 4) Widen column G to fit the new, much longer formula results.
$ws.Columns("G").ColumnWidth = 101.6

# 5) Update the view: scroll/zoom plus reselect the newly filled G4:G36 range.
$excel.ActiveWindow.Zoom = 130
$ws.Range("G4:G36").Select()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 2
